$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for bf6fe2c5... (row3) and c4986af8... (row4)
# This shared string is also used by de-de's "Correspond Handoff Datetime" (H3/H4),
# so those cells move together with the Overview cells.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-02 14:19:32"
$wsOverview.Range("G4").Value = "2016-09-02 14:19:32"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority: ht -> mt
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H3").Value = "2016-09-02 14:19:27"
$wsZhCn.Range("H4").Value = "2016-09-02 14:19:27"
# Correspond Handback DateTime
$wsZhCn.Range("K3").Value = "2016-09-02 14:19:45"
$wsZhCn.Range("K4").Value = "2016-09-02 14:19:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority: ht -> mt
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime (shares its text with Overview's G3/G4 above)
$wsDeDe.Range("H3").Value = "2016-09-02 14:19:32"
$wsDeDe.Range("H4").Value = "2016-09-02 14:19:32"
# Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-09-02 14:19:52"
$wsDeDe.Range("K4").Value = "2016-09-02 14:19:52"
